$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve column D price strings as text (they look numeric, e.g. "27.108.35")
# by setting the Text number format before assigning, matching the source data
# which stores prices as literal strings, not numbers.
$dCells = @("D2","D3","D4","D5","D7","D8","D9","D10","D11","D12","D13","D14","D15","D16","D18","D20","D22","D23","D24","D25","D27","D28","D29","D30","D31","D32","D33","D34","D35","D37","D38","D40","D42","D43","D44","D45","D46","D48","D49","D50","D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.108.35"
$ws.Range("E2").Value = "  -1.03%  "
$ws.Range("D3").Value = "1.824.68"
$ws.Range("E3").Value = "  -0.96%  "
$ws.Range("D4").Value = "1.009"
$ws.Range("E4").Value = "  -0.61%  "
$ws.Range("D5").Value = "311.39"
$ws.Range("E5").Value = "  -1.31%  "
$ws.Range("D7").Value = "0.4651"
$ws.Range("E7").Value = "  -2.04%  "
$ws.Range("D8").Value = "0.3639"
$ws.Range("E8").Value = "  -1.71%  "
$ws.Range("D9").Value = "0.07298"
$ws.Range("E9").Value = "  -2.35%  "
$ws.Range("D10").Value = "0.8703"
$ws.Range("E10").Value = "  -1.81%  "
$ws.Range("D11").Value = "20.19"
$ws.Range("E11").Value = "  -1.63%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "0.07628"
$ws.Range("E12").Value = "  +3.42%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.853.93"
$ws.Range("E13").Value = "  -0.26%  "
$ws.Range("B14").Value = "Litecoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D14").Value = "92.77"
$ws.Range("E14").Value = "  -0.61%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "5.345"
$ws.Range("E15").Value = "  -2.68%  "
$ws.Range("D16").Value = "6.483"
$ws.Range("E16").Value = "  -1.53%  "
$ws.Range("E17").Value = "  -0.59%  "
$ws.Range("D18").Value = "0.000008652"
$ws.Range("E18").Value = "  -2.41%  "
$ws.Range("E19").Value = "  -0.46%  "
$ws.Range("D20").Value = "27.245.33"
$ws.Range("E20").Value = "  -0.62%  "
$ws.Range("E21").Value = "  -2.43%  "
$ws.Range("D22").Value = "5.196"
$ws.Range("E22").Value = "  -2.91%  "
$ws.Range("D23").Value = "10.56"
$ws.Range("E23").Value = "  -1.57%  "
$ws.Range("D24").Value = "2.083.36"
$ws.Range("D25").Value = "151.69"
$ws.Range("E25").Value = "  -0.32%  "
$ws.Range("E26").Value = "  -2.49%  "
$ws.Range("D27").Value = "18.28"
$ws.Range("E27").Value = "  -2.02%  "
$ws.Range("D28").Value = "2.111"
$ws.Range("E28").Value = "  -3.29%  "
$ws.Range("D29").Value = "116.14"
$ws.Range("E29").Value = "  -1.68%  "
$ws.Range("D30").Value = "5.092"
$ws.Range("E30").Value = "  -3.55%  "
$ws.Range("D31").Value = "0.08921"
$ws.Range("E31").Value = "  -0.56%  "
$ws.Range("D32").Value = "2.959"
$ws.Range("E32").Value = "  +0.38%  "
$ws.Range("D33").Value = "0.7355"
$ws.Range("E33").Value = "  -3.51%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "4.451"
$ws.Range("E34").Value = "  -2.54%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "1.142"
$ws.Range("E35").Value = "  -3.27%  "
$ws.Range("E36").Value = "  -0.41%  "
$ws.Range("D37").Value = "2.551"
$ws.Range("E37").Value = "  +7.34%  "
$ws.Range("D38").Value = "0.05265"
$ws.Range("E38").Value = "  -1.96%  "
$ws.Range("E39").Value = "  -3.52%  "
$ws.Range("D40").Value = "0.01919"
$ws.Range("E40").Value = "  -2.25%  "
$ws.Range("E41").Value = "  -2.31%  "
$ws.Range("D42").Value = "7.129"
$ws.Range("E42").Value = "  -2.73%  "
$ws.Range("D43").Value = "0.5223"
$ws.Range("E43").Value = "  -2.56%  "
$ws.Range("D44").Value = "0.1634"
$ws.Range("E44").Value = "  -2.19%  "
$ws.Range("D45").Value = "8.275"
$ws.Range("E45").Value = "  -3.23%  "
$ws.Range("D46").Value = "0.4876"
$ws.Range("E46").Value = "  -2.13%  "
$ws.Range("E47").Value = "  -0.49%  "
$ws.Range("D48").Value = "103.89"
$ws.Range("E48").Value = "  -1.30%  "
$ws.Range("D49").Value = "10.13"
$ws.Range("E49").Value = "  -3.66%  "
$ws.Range("D50").Value = "1.639"
$ws.Range("E50").Value = "  -2.61%  "
$ws.Range("D51").Value = "0.06245"
$ws.Range("E51").Value = "  -1.43%  "
